# spring 24 week 15 inputs
# Append the week's matchup rows (1619-1635) to the "Nine" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Nine")

$data = @(
  @(6, 13, 7, 7),
  @(4, 13, 2, 7),
  @(7, 15, 4, 5),
  @(5, 6, 9, 14),
  @(5, 13, 4, 7),
  @(4, 6, 5, 14),
  @(5, 15, 4, 5),
  @(5, 12, 4, 8),
  @(2, 16, 4, 4),
  @(6, 13, 5, 7),
  @(3, 4, 5, 16),
  @(3, 15, 2, 5),
  @(3, 15, 4, 5),
  @(8, 6, 7, 14),
  @(2, 3, 3, 17),
  @(5, 16, 4, 4),
  @(8, 18, 6, 2)
)

$startRow = 1619
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
    $ws.Cells.Item($row, 4).Value = $data[$i][3]
}

# Match the author's final view state: scrolled down near the new rows,
# with C1619 (first blank cell of the next row) selected.
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 1616
$ws.Range("C1619").Select()
